$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 10 (Number = DEI-3-2): add a note in column G
$ws.Range("G10").Value = "Check defects, check defects and downtime. Check states"

# Row 11 (Number = DEI-3-3): add real hours estimate in column B
$ws.Range("B11").Value = "3 h"

# Update the active selection to B13
$ws.Range("B13").Select()
